$wb = $excel.ActiveWorkbook

# Add a new worksheet right after Sheet1 and populate it with the
# address / contact-details table ("Sheet2" in the target workbook).
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "City"
$ws2.Range("B1").Value = "Address 1"
$ws2.Range("C1").Value = "Address 2"
$ws2.Range("D1").Value = "Postal Code"
$ws2.Range("E1").Value = "Phone number"
$ws2.Range("F1").Value = "FaxNumber"

$ws2.Range("A2").Value = "Delhi"
$ws2.Range("B2").Value = "Random #103 Delhi India"
$ws2.Range("C2").Value = "Random #103 Delhi India"
$ws2.Range("D2").Value = 123456
$ws2.Range("E2").Value = 9876543210
$ws2.Range("F2").Value = 55555555

# Column widths for the new sheet (matches the author's autofit widths as
# closely as this host's column-width model allows).
$ws2.Columns.Item(1).ColumnWidth = 13.666666666666666
$ws2.Columns.Item(2).ColumnWidth = 25.833333333333332
$ws2.Columns.Item(3).ColumnWidth = 16.166666666666668
$ws2.Columns.Item(4).ColumnWidth = 17
$ws2.Columns.Item(5).ColumnWidth = 16.666666666666668
$ws2.Columns.Item(6).ColumnWidth = 28.333333333333332

$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# Selection / active-sheet state matches the authored workbook: Sheet2
# becomes the visible (tab-selected) sheet with F9 selected.
$ws2.Range("F9").Select()
$ws2.Activate()
